$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between column B and C for the rows that changed
$rows = @(2, 5, 7, 10, 12, 13, 16)
foreach ($r in $rows) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value = $cVal
    $ws.Cells.Item($r, 3).Value = $bVal
}

# Update the active selection to B19
$ws.Range("B19").Select()
